$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at row 757, shifting existing rows 757:786 down to 760:789
$ws.Rows("757:759").Insert()

# New data rows (Terminal La Palmera de La Serena - Mandarina - Clementina, 2023-05-29)
$data = @(
    @(8, 'Terminal La Palmera de La Serena', 'Coquimbo', '2023-05-29', 4, 'Fruta', 100102, 'Cítricos', 100102004, 'Mandarina', 'Clementina', 'Especial', 20, 290000, 300000, 295000, '$/bins (450 kilos)', 'Provincia de Limarí', 656, 450),
    @(8, 'Terminal La Palmera de La Serena', 'Coquimbo', '2023-05-29', 4, 'Fruta', 100102, 'Cítricos', 100102004, 'Mandarina', 'Clementina', 'Primera', 20, 260000, 270000, 265000, '$/bins (450 kilos)', 'Provincia de Limarí', 589, 450),
    @(8, 'Terminal La Palmera de La Serena', 'Coquimbo', '2023-05-29', 4, 'Fruta', 100102, 'Cítricos', 100102004, 'Mandarina', 'Clementina', 'Segunda', 20, 220000, 230000, 225000, '$/bins (450 kilos)', 'Provincia de Limarí', 500, 450)
)

$startRow = 757
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    for ($c = 1; $c -le $row.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
